# Updated symbol list on Sun Jan  8 19:38:39 UTC 2023 with GitHub Actions
#
# The "Price" (column D) and "Volume(1h)" (column E) columns on the active
# sheet are refreshed with newer quotes for the rows below. Every value in
# this sheet is stored as literal text (e.g. "265.79", "1.61%"), so each
# new value is written with a leading apostrophe to force Excel to keep it
# as text instead of re-interpreting it as a number/percentage, and the
# cell style is reset back to "Normal" right after so no stray
# number-format/quote-prefix style is left attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "265.79"
    "E2"  = "1.61%"
    "D3"  = "26.72"
    "E3"  = "-1.76%"
    "D4"  = "4.694"
    "E4"  = "-0.33%"
    "D5"  = "0.06088"
    "E5"  = "-1.75%"
    "D6"  = "6.733"
    "E6"  = "0.29%"
    "D7"  = "0.8498"
    "E7"  = "-0.18%"
    "E8"  = "-1.13%"
    "E9"  = "-0.12%"
    "D10" = "0.05021"
    "E10" = "10.17%"
    "D11" = "0.07085"
    "E11" = "-0.04%"
    "D12" = "0.03183"
    "E12" = "1.75%"
    "D13" = "0.09013"
    "E13" = "-0.45%"
    "D14" = "0.001533"
    "E14" = "0.19%"
    "D15" = "0.0006034"
    "E15" = "-1.90%"
    "D16" = "0.005958"
    "E16" = "-1.92%"
    "D17" = "3.460"
    "E17" = "0.01%"
    "E18" = "0.09%"
    "D19" = "2.263"
    "E19" = "3.18%"
    "E21" = "-0.66%"
    "D22" = "4.083"
    "E22" = "0.06%"
    "D23" = "0.04222"
    "E23" = "-0.21%"
    "D24" = "0.001186"
    "E24" = "-2.47%"
    "D25" = "0.004139"
    "E25" = "8.85%"
    "E26" = "0.04%"
    "D27" = "0.0001682"
    "E27" = "5.07%"
    "D40" = "0.03912"
    "E40" = "-0.27%"
    "E41" = "-0.02%"
    "D42" = "0.004183"
    "E42" = "1.19%"
    "D43" = "0.002121"
    "E43" = "-2.86%"
    "D44" = "0.01261"
    "E44" = "-8.90%"
    "D45" = "0.00005120"
    "E45" = "-0.27%"
    "E46" = "0.04%"
    "E47" = "-31.79%"
    "D48" = "0.1520"
    "E48" = "-9.38%"
    "E49" = "0.04%"
    "E50" = "0.04%"
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $updates[$cellRef]
    $cell.Style = "Normal"
}
